$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Update VENTA/POR CUMPLIR/CUMPLIMIENTO for rows whose venta resets to 0 ---
# Row 3 - 240X80 PORCELANATO
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 5504.61890386263
$ws.Range("F3").Value = 0

# Row 5 - GRIFERIAS
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 150
$ws.Range("F5").Value = 0

# Row 6 - INODOROS
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 2907.58368146026
$ws.Range("F6").Value = 0

# Row 7 - LAVABOS
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 886.711016287574
$ws.Range("F7").Value = 0

# Row 11 - PIEDRA SINTERIZADA
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 5844.44916370549
$ws.Range("F11").Value = 0

# Row 12 - PORCELANATO (presupuesto also updated)
$ws.Range("C12").Value = 37739.74
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 37739.74
$ws.Range("F12").Value = 0

# --- Remove the "SAL SOLUBLE" row (row 14); TOTAL row shifts up from 15 to 14 ---
$ws.Rows.Item(14).Delete()

# --- Refresh the TOTAL row (now row 14) with the recalculated sums ---
$ws.Range("C14").Value = 55424.74147880389
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 55424.74147880389
$ws.Range("F14").Value = 0

# --- Resize columns D, E, F (character widths, offset-compensated to match target) ---
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668
